# Fig 5 schematic edits:
#   - "1-day of 15% PEG 3350 in drinking water (N = 18)" -> "1-day of 15% PEG in drinking water (N = 18)"
#   - "3-day recovery + 1-day PEG 3350 + FMT (N = 6)"     -> "3-day recovery + 1-day PEG + FMT (N = 6)"
#   - "3-day recovery + 1-day PEG 3350 + PBS (N = 12)"    -> "3-day recovery + 1-day PEG + PBS (N = 12)"
#
# The three legend lines live as separate runs inside one grouped textbox
# ("Google Shape;698;p64", inside group "Group 6") on slide 1. We locate the
# shape by scanning (rather than trusting a hard-coded index chain), then
# edit each affected run's .Text directly so we only touch the characters
# that changed and leave every other run's formatting (color, font, etc.)
# untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Type -eq 6) {
        # msoGroup - only need to look one level deep for this deck
        $items = $shp.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            $inner = $items.Item($j)
            if ($inner.HasTextFrame) {
                if ($inner.TextFrame.HasText) {
                    $t = $inner.TextFrame.TextRange.Text
                    if ($t -like "*PEG 3350*") {
                        $targetShape = $inner
                    }
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count

    for ($i = 1; $i -le $paraCount; $i++) {
        $para = $tr.Paragraphs($i)
        $runCount = $para.Runs().Count
        for ($k = 1; $k -le $runCount; $k++) {
            $run = $para.Runs($k)
            $txt = $run.Text
            # The last run of a paragraph reports a trailing CR (paragraph
            # mark) as part of .Text; strip it before editing so we don't
            # bake a literal newline into the run's <a:t> when writing back.
            $trimmed = $txt.TrimEnd("`r", "`n")
            if ($trimmed -like "*PEG 3350*") {
                $run.Text = $trimmed -replace "PEG 3350", "PEG"
            }
        }
    }
}
